$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# spawn rule event update: waitGenerateTime for row 14 drops from 100 to 5
$ws.Range("C14").Value = 5

# Update the last-saved cell selection/cursor position recorded in the sheet view
$ws.Range("J21").Select()
